$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("INCIO SANCHEZ PAOLA KATHERINE", 80),
    @("GUEVARA IDROGO DENNIS PERCY", 76),
    @("TANTALEAN BUSTAMANTE ESTALIN YOEL", 74),
    @("HUAYHUA VALDIVIA LUZ EXMILDA", 69),
    @("DELGADO VASQUEZ FLOR MAGALY", 68),
    @("LINARES PEREZ YANASELY", 67),
    @("MEDINA TAPIA ANA YULI", 66),
    @("LOZADA ROJAS LUZ ELENA", 66),
    @("CAMPOS PEREZ YOVERLY", 65),
    @("CHAVEZ VILLANUEVA SILVIA JANETH", 65),
    @("VASQUEZ SILVA ALOIS ADOLF", 65),
    @("PEREZ LINARES TATHIANA", 65),
    @("MONDRAGON HERNANDEZ WILMER JUNIOR", 64),
    @("SOTO LOZANO LUZDINA", 62)
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $row = $row + 1
}
